# Add Job API - adds a new applicant (Aman Kumar) to the "Users" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$ws.Range("A3").Value = "U#00002"
$ws.Range("B3").Value = "Aman"
$ws.Range("C3").Value = "Kumar"

# Mobile number looks numeric, force it to stay text (matches existing D2 style)
$ws.Range("D3").Value = "'9031398069"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "aman@gmail.com"

# Password also looks numeric, force it to stay text (matches existing F2 style)
$ws.Range("F3").Value = "'1234"
$ws.Range("F3").Style = "Normal"

$ws.Range("G3").Value = "java,django,cpp,bigData,networking"
